$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.01997857905520561
$ws.Cells.Item(2, 4).Value = 0.01930819441224685
$ws.Cells.Item(2, 5).Value = 0.130636675233994
$ws.Cells.Item(2, 6).Value = 0.4424814950166223
$ws.Cells.Item(2, 7).Value = 0.289050929862789
$ws.Cells.Item(2, 8).Value = 0.4517458746030272
$ws.Cells.Item(2, 9).Value = 0.3249404618561798
$ws.Cells.Item(2, 11).Value = 0.7499470307241438
$ws.Cells.Item(2, 13).Value = 0.3056097893511591
$ws.Cells.Item(2, 14).Value = 0.993287885816855
$ws.Cells.Item(2, 15).Value = 1.406571826920398
$ws.Cells.Item(3, 3).Value = 0.01740730519347267
$ws.Cells.Item(3, 4).Value = 0.01708077576778066
$ws.Cells.Item(3, 5).Value = 0.1244272219009446
$ws.Cells.Item(3, 6).Value = 0.4384638198585122
$ws.Cells.Item(3, 7).Value = 0.286334960648297
$ws.Cells.Item(3, 8).Value = 0.4540420131334386
$ws.Cells.Item(3, 9).Value = 0.3226345430633053
$ws.Cells.Item(3, 11).Value = 0.6543143975407588
$ws.Cells.Item(3, 13).Value = 0.2715890736685154
$ws.Cells.Item(3, 14).Value = 0.994344804798672
$ws.Cells.Item(3, 15).Value = 1.405364914913903
$ws.Cells.Item(4, 3).Value = 0.0158244275320456
$ws.Cells.Item(4, 4).Value = 0.01570547196724448
$ws.Cells.Item(4, 5).Value = 0.1207400044384528
$ws.Cells.Item(4, 6).Value = 0.4363337317295688
$ws.Cells.Item(4, 7).Value = 0.2849324706649341
$ws.Cells.Item(4, 8).Value = 0.4556960244562021
$ws.Cells.Item(4, 9).Value = 0.3214676974139756
$ws.Cells.Item(4, 11).Value = 0.5953679209982852
$ws.Cells.Item(4, 13).Value = 0.2507364586208638
$ws.Cells.Item(4, 14).Value = 0.9953708979673621
$ws.Cells.Item(4, 15).Value = 1.405688021901184
$ws.Cells.Item(5, 3).Value = 0.01517839224198525
$ws.Cells.Item(5, 4).Value = 0.01514313469117212
$ws.Cells.Item(5, 5).Value = 0.1192687151943304
$ws.Cells.Item(5, 6).Value = 0.4355503156784621
$ws.Cells.Item(5, 7).Value = 0.2844274486394482
$ws.Cells.Item(5, 8).Value = 0.4564314478919869
$ws.Cells.Item(5, 9).Value = 0.3210547459319493
$ws.Cells.Item(5, 11).Value = 0.5712911943524546
$ws.Cells.Item(5, 13).Value = 0.2422481540065604
$ws.Cells.Item(5, 14).Value = 0.9958841124018747
$ws.Cells.Item(5, 15).Value = 1.406086929509385
$ws.Cells.Item(6, 3).Value = 0.01507105914376439
$ws.Cells.Item(6, 4).Value = 0.01504964594235503
$ws.Cells.Item(6, 5).Value = 0.1190262906494155
$ws.Cells.Item(6, 6).Value = 0.4354253382323421
$ws.Cells.Item(6, 7).Value = 0.2843476019807198
$ws.Cells.Item(6, 8).Value = 0.4565572728635772
$ws.Cells.Item(6, 9).Value = 0.3209899514690555
$ws.Cells.Item(6, 11).Value = 0.5672899663153714
$ws.Cells.Item(6, 13).Value = 0.2408392452042065
$ws.Cells.Item(6, 14).Value = 0.9959750805493925
$ws.Cells.Item(6, 15).Value = 1.40616929807193
$ws.Cells.Item(7, 3).Value = 0.01581571887872002
$ws.Cells.Item(7, 4).Value = 0.01569789569422397
$ws.Cells.Item(7, 5).Value = 0.1207200357942071
$ws.Cells.Item(7, 6).Value = 0.4363228238143506
$ws.Cells.Item(7, 7).Value = 0.2849253907068388
$ws.Cells.Item(7, 8).Value = 0.4557056940249424
$ws.Cells.Item(7, 9).Value = 0.3214618750396525
$ws.Cells.Item(7, 11).Value = 0.5950434365608999
$ws.Cells.Item(7, 13).Value = 0.2506219444437505
$ws.Cells.Item(7, 14).Value = 0.9953774340907842
$ws.Cells.Item(7, 15).Value = 1.405692320137604
$ws.Cells.Item(8, 3).Value = 0.01909287458900621
$ws.Cells.Item(8, 4).Value = 0.01854179137644252
$ws.Cells.Item(8, 5).Value = 0.1284694715171284
$ws.Cells.Item(8, 6).Value = 0.4410262547510655
$ws.Cells.Item(8, 7).Value = 0.2880593364673203
$ws.Cells.Item(8, 8).Value = 0.4524869051496836
$ws.Cells.Item(8, 9).Value = 0.324093660860818
$ws.Cells.Item(8, 11).Value = 0.7170211303663336
$ws.Cells.Item(8, 13).Value = 0.293872029196045
$ws.Cells.Item(8, 14).Value = 0.9935741753677974
$ws.Cells.Item(8, 15).Value = 1.4059346019299
$ws.Cells.Item(9, 3).Value = 0.02548566431050858
$ws.Cells.Item(9, 4).Value = 0.02405653896126125
$ws.Cells.Item(9, 5).Value = 0.144673522259545
$ws.Cells.Item(9, 6).Value = 0.4529264250736205
$ws.Cells.Item(9, 7).Value = 0.2963168538882286
$ws.Cells.Item(9, 8).Value = 0.4481125409603948
$ws.Cells.Item(9, 9).Value = 0.3312340079804557
$ws.Cells.Item(9, 11).Value = 0.9543497121239568
$ws.Cells.Item(9, 13).Value = 0.3789687357264455
$ws.Cells.Item(9, 14).Value = 0.9930209894473592
$ws.Cells.Item(9, 15).Value = 1.414871548429517
$ws.Cells.Item(10, 3).Value = 0.03016081821421324
$ws.Cells.Item(10, 4).Value = 0.02806890730363421
$ws.Cells.Item(10, 5).Value = 0.1572111495284361
$ws.Cells.Item(10, 6).Value = 0.4633093500992373
$ws.Cells.Item(10, 7).Value = 0.3036836032238739
$ws.Cells.Item(10, 8).Value = 0.4460807828923805
$ws.Cells.Item(10, 9).Value = 0.3376930674264997
$ws.Cells.Item(10, 11).Value = 1.127507153572765
$ws.Cells.Item(10, 13).Value = 0.4416628493032846
$ws.Cells.Item(10, 14).Value = 0.9944218769293371
$ws.Cells.Item(10, 15).Value = 1.426624728478913
$ws.Cells.Item(11, 3).Value = 0.03228276463171653
$ws.Cells.Item(11, 4).Value = 0.02988541477677842
$ws.Cells.Item(11, 5).Value = 0.1630562042624675
$ws.Cells.Item(11, 6).Value = 0.4683907127435774
$ws.Cells.Item(11, 7).Value = 0.3073199398837403
$ws.Cells.Item(11, 8).Value = 0.4454133764491672
$ws.Cells.Item(11, 9).Value = 0.3408962843025662
$ws.Cells.Item(11, 11).Value = 1.206005219862504
$ws.Cells.Item(11, 13).Value = 0.4702222675767018
$ws.Cells.Item(11, 14).Value = 0.9954494312490851
$ws.Cells.Item(11, 15).Value = 1.433104246168028
$ws.Cells.Item(12, 3).Value = 0.03308557178243632
$ws.Cells.Item(12, 4).Value = 0.0305719896519463
$ws.Cells.Item(12, 5).Value = 0.1652902362420292
$ws.Cells.Item(12, 6).Value = 0.4703664974995121
$ws.Cells.Item(12, 7).Value = 0.3087381268250482
$ws.Cells.Item(12, 8).Value = 0.4451975927627245
$ws.Cells.Item(12, 9).Value = 0.3421474490558012
$ws.Cells.Item(12, 11).Value = 1.2356897666848
$ws.Cells.Item(12, 13).Value = 0.4810425609517921
$ws.Cells.Item(12, 14).Value = 0.9958944681702633
$ws.Cells.Item(12, 15).Value = 1.435721227412188
$ws.Cells.Item(13, 3).Value = 0.03291270586255735
$ws.Cells.Item(13, 4).Value = 0.03042418180012874
$ws.Cells.Item(13, 5).Value = 0.164808175348746
$ws.Cells.Item(13, 6).Value = 0.4699386813800075
$ws.Cells.Item(13, 7).Value = 0.3084308603418435
$ws.Cells.Item(13, 8).Value = 0.4452424221199465
$ws.Cells.Item(13, 9).Value = 0.3418762892623377
$ws.Cells.Item(13, 11).Value = 1.229298519619533
$ws.Cells.Item(13, 13).Value = 0.4787119778719102
$ws.Cells.Item(13, 14).Value = 0.9957961373478099
$ws.Cells.Item(13, 15).Value = 1.435150343050509
$ws.Cells.Item(14, 3).Value = 0.03234882685968898
$ws.Cells.Item(14, 4).Value = 0.02994192593241252
$ws.Cells.Item(14, 5).Value = 0.163239584277072
$ws.Cells.Item(14, 6).Value = 0.4685522274571454
$ws.Cells.Item(14, 7).Value = 0.3074357885466128
$ws.Cells.Item(14, 8).Value = 0.4453948832105823
$ws.Cells.Item(14, 9).Value = 0.3409984528475078
$ws.Cells.Item(14, 11).Value = 1.208448216941406
$ws.Cells.Item(14, 13).Value = 0.4711123514631907
$ws.Cells.Item(14, 14).Value = 0.9954849250106292
$ws.Cells.Item(14, 15).Value = 1.433316271087051
$ws.Cells.Item(15, 3).Value = 0.03200333879375705
$ws.Cells.Item(15, 4).Value = 0.02964636017001965
$ws.Cells.Item(15, 5).Value = 0.162281472493504
$ws.Cells.Item(15, 6).Value = 0.4677097045524121
$ws.Cells.Item(15, 7).Value = 0.3068316471423316
$ws.Cells.Item(15, 8).Value = 0.4454930822450507
$ws.Cells.Item(15, 9).Value = 0.3404657265610211
$ws.Cells.Item(15, 11).Value = 1.195671414946048
$ws.Cells.Item(15, 13).Value = 0.4664580653917341
$ws.Cells.Item(15, 14).Value = 0.9953015756942989
$ws.Cells.Item(15, 15).Value = 1.432214131458551
$ws.Cells.Item(16, 3).Value = 0.03002204418491772
$ws.Cells.Item(16, 4).Value = 0.02795001496202332
$ws.Cells.Item(16, 5).Value = 0.1568320344897387
$ws.Cells.Item(16, 6).Value = 0.4629844847447657
$ws.Cells.Item(16, 7).Value = 0.3034517136415928
$ws.Cells.Item(16, 8).Value = 0.44612956914915
$ws.Cells.Item(16, 9).Value = 0.3374890677767723
$ws.Cells.Item(16, 11).Value = 1.122371487442138
$ws.Cells.Item(16, 13).Value = 0.4397972007193403
$ws.Cells.Item(16, 14).Value = 0.9943625603009849
$ws.Cells.Item(16, 15).Value = 1.426224109783021
$ws.Cells.Item(17, 3).Value = 0.02880532725588125
$ws.Cells.Item(17, 4).Value = 0.02690709509240463
$ws.Cells.Item(17, 5).Value = 0.15352544728308
$ws.Cells.Item(17, 6).Value = 0.4601774998119694
$ws.Cells.Item(17, 7).Value = 0.301451399276857
$ws.Cells.Item(17, 8).Value = 0.4465858287261995
$ws.Cells.Item(17, 9).Value = 0.335730900186654
$ws.Cells.Item(17, 11).Value = 1.077333350877609
$ws.Cells.Item(17, 13).Value = 0.4234515931628096
$ws.Cells.Item(17, 14).Value = 0.993886323236481
$ws.Cells.Item(17, 15).Value = 1.422839893253979
$ws.Cells.Item(18, 3).Value = 0.02810505385315309
$ws.Cells.Item(18, 4).Value = 0.02630641502817355
$ws.Cells.Item(18, 5).Value = 0.1516369000140045
$ws.Cells.Item(18, 6).Value = 0.4585966990360646
$ws.Cells.Item(18, 7).Value = 0.3003277018929253
$ws.Cells.Item(18, 8).Value = 0.446872432327595
$ws.Cells.Item(18, 9).Value = 0.334744583388499
$ws.Cells.Item(18, 11).Value = 1.051403101958158
$ws.Cells.Item(18, 13).Value = 0.414053768323555
$ws.Cells.Item(18, 14).Value = 0.9936491510915744
$ws.Cells.Item(18, 15).Value = 1.420999994316674
$ws.Cells.Item(19, 3).Value = 0.02786787714758532
$ws.Cells.Item(19, 4).Value = 0.02610289545086886
$ws.Cells.Item(19, 5).Value = 0.1509997488302375
$ws.Cells.Item(19, 6).Value = 0.4580672528803191
$ws.Cells.Item(19, 7).Value = 0.2999518399656296
$ws.Cells.Item(19, 8).Value = 0.4469736230663841
$ws.Cells.Item(19, 9).Value = 0.3344149137124575
$ws.Cells.Item(19, 11).Value = 1.042619240920942
$ws.Cells.Item(19, 13).Value = 0.4108724749048775
$ws.Cells.Item(19, 14).Value = 0.993575166100527
$ws.Cells.Item(19, 15).Value = 1.420395334363775
$ws.Cells.Item(20, 3).Value = 0.02893489569001417
$ws.Cells.Item(20, 4).Value = 0.02701820086718243
$ws.Cells.Item(20, 5).Value = 0.153876059199753
$ws.Cells.Item(20, 6).Value = 0.4604728194165872
$ws.Cells.Item(20, 7).Value = 0.3016615579968516
$ws.Cells.Item(20, 8).Value = 0.4465347569750548
$ws.Cells.Item(20, 9).Value = 0.3359154789308434
$ws.Cells.Item(20, 11).Value = 1.082130390284362
$ws.Cells.Item(20, 13).Value = 0.4251912259359329
$ws.Cells.Item(20, 14).Value = 0.9939332177765863
$ws.Cells.Item(20, 15).Value = 1.423189111837161
$ws.Cells.Item(21, 3).Value = 0.03251447184210576
$ws.Cells.Item(21, 4).Value = 0.03008361169507623
$ws.Cells.Item(21, 5).Value = 0.1636997552663857
$ws.Cells.Item(21, 6).Value = 0.468958061956755
$ws.Cells.Item(21, 7).Value = 0.3077269460834628
$ws.Cells.Item(21, 8).Value = 0.4453490988068154
$ws.Cells.Item(21, 9).Value = 0.3412552579998405
$ws.Cells.Item(21, 11).Value = 1.214573582707487
$ws.Cells.Item(21, 13).Value = 0.4733443995354065
$ws.Cells.Item(21, 14).Value = 0.9955748193963956
$ws.Cells.Item(21, 15).Value = 1.433850546514236
$ws.Cells.Item(22, 3).Value = 0.03484967173176301
$ws.Cells.Item(22, 4).Value = 0.03207945572977877
$ws.Cells.Item(22, 5).Value = 0.170240520242956
$ws.Cells.Item(22, 6).Value = 0.4748043624614979
$ws.Cells.Item(22, 7).Value = 0.3119311683484369
$ws.Cells.Item(22, 8).Value = 0.4447895739946688
$ws.Cells.Item(22, 9).Value = 0.3449676656562701
$ws.Cells.Item(22, 11).Value = 1.300893227952145
$ws.Cells.Item(22, 13).Value = 0.5048470109120586
$ws.Cells.Item(22, 14).Value = 0.9969735831030988
$ws.Cells.Item(22, 15).Value = 1.441770603369235
$ws.Cells.Item(23, 3).Value = 0.03360373365879354
$ws.Cells.Item(23, 4).Value = 0.03101494305730768
$ws.Cells.Item(23, 5).Value = 0.1667384836396053
$ws.Cells.Item(23, 6).Value = 0.4716565394075616
$ws.Cells.Item(23, 7).Value = 0.3096652642746136
$ws.Cells.Item(23, 8).Value = 0.4450684918797805
$ws.Cells.Item(23, 9).Value = 0.3429658955706429
$ws.Cells.Item(23, 11).Value = 1.254845341091141
$ws.Cells.Item(23, 13).Value = 0.488030645859979
$ws.Cells.Item(23, 14).Value = 0.9961972823538048
$ws.Cells.Item(23, 15).Value = 1.437456258540976
$ws.Cells.Item(24, 3).Value = 0.02887632018679653
$ws.Cells.Item(24, 4).Value = 0.02696797334355949
$ws.Cells.Item(24, 5).Value = 0.1537175087821723
$ws.Cells.Item(24, 6).Value = 0.4603392027239224
$ws.Cells.Item(24, 7).Value = 0.3015664632995936
$ws.Cells.Item(24, 8).Value = 0.4465577708403856
$ws.Cells.Item(24, 9).Value = 0.3358319546417761
$ws.Cells.Item(24, 11).Value = 1.079961764578798
$ws.Cells.Item(24, 13).Value = 0.4244047395566213
$ws.Cells.Item(24, 14).Value = 0.993911902660372
$ws.Cells.Item(24, 15).Value = 1.423030900797357
$ws.Cells.Item(25, 3).Value = 0.02375996493289279
$ws.Cells.Item(25, 4).Value = 0.02257145367626379
$ws.Cells.Item(25, 5).Value = 0.1401801476286693
$ws.Cells.Item(25, 6).Value = 0.449419711340532
$ws.Cells.Item(25, 7).Value = 0.2938555899724804
$ws.Cells.Item(25, 8).Value = 0.4490883953416898
$ws.Cells.Item(25, 9).Value = 0.3290897853869268
$ws.Cells.Item(25, 11).Value = 0.8903532428073504
$ws.Cells.Item(25, 13).Value = 0.355917398774146
$ws.Cells.Item(25, 14).Value = 0.9928524699872128
$ws.Cells.Item(25, 15).Value = 1.411544998278345

Write-Output "done"